$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ALC row 41
$ws1.Range("H41").Value = 579.0833
$ws1.Range("I41").Value = 781.25
$ws1.Range("J41").Value = 174.75
$ws1.Range("K41").Value = 781.25
$ws1.Range("L41").Value = 174.75
$ws1.Range("M41").Value = -341.25
$ws1.Range("N41").Value = -1054.75

# ALC row 43
$ws1.Range("H43").Value = 3000
$ws1.Range("J43").Value = 3000
$ws1.Range("L43").Value = 3000
$ws1.Range("N43").Value = -3138

# ALC row 76
$ws1.Range("H76").Value = 10870.2
$ws1.Range("I76").Value = 9936.77
$ws1.Range("J76").Value = 16937.5
$ws1.Range("K76").Value = 9936.77
$ws1.Range("L76").Value = 16937.5
$ws1.Range("M76").Value = -9621.77
$ws1.Range("N76").Value = -17567.5

# ALC row 79
$ws1.Range("H79").Value = 10870.2
$ws1.Range("I79").Value = 9936.77
$ws1.Range("J79").Value = 16937.5
$ws1.Range("K79").Value = 9936.77
$ws1.Range("L79").Value = 16937.5
$ws1.Range("M79").Value = -8844.77
$ws1.Range("N79").Value = -19121.5

# ALC row 112
$ws1.Range("H112").Value = 4001.4119
$ws1.Range("I112").Value = 1950
$ws1.Range("J112").Value = 4274.933
$ws1.Range("K112").Value = 5850
$ws1.Range("L112").Value = 12824.799
$ws1.Range("M112").Value = -4742
$ws1.Range("N112").Value = -15040.799

# ALC row 135
$ws1.Range("H135").Value = 1166.0322
$ws1.Range("I135").Value = 778.1429000000001
$ws1.Range("K135").Value = 7003.2861
$ws1.Range("M135").Value = -4468.2861

# ALC row 137
$ws1.Range("H137").Value = 1975.6666
$ws1.Range("I137").Value = 1971.2858
$ws1.Range("J137").Value = 1991
$ws1.Range("K137").Value = 5913.857400000001
$ws1.Range("L137").Value = 5973
$ws1.Range("M137").Value = -3363.857400000001
$ws1.Range("N137").Value = -11073

# ALC row 138
$ws1.Range("H138").Value = 5354.431
$ws1.Range("J138").Value = 7693.027
$ws1.Range("L138").Value = 23079.081
$ws1.Range("N138").Value = -33359.081

# ARM row 32
$ws2.Range("H32").Value = 8073.2295
$ws2.Range("I32").Value = 640.4583
$ws2.Range("K32").Value = 640.4583
$ws2.Range("M32").Value = -353.4583

# ARM row 74
$ws2.Range("H74").Value = 2283.3333
$ws2.Range("I74").Value = 2273.2559
$ws2.Range("J74").Value = 2500
$ws2.Range("K74").Value = 2273.2559
$ws2.Range("L74").Value = 2500
$ws2.Range("M74").Value = -1399.2559
$ws2.Range("N74").Value = -4248

# ARM row 77
$ws2.Range("H77").Value = 2283.3333
$ws2.Range("I77").Value = 2273.2559
$ws2.Range("J77").Value = 2500
$ws2.Range("K77").Value = 11366.2795
$ws2.Range("L77").Value = 12500
$ws2.Range("M77").Value = -6998.279500000001
$ws2.Range("N77").Value = -21236

# ARM row 124
$ws2.Range("H124").Value = 25000
$ws2.Range("J124").Value = 25000
$ws2.Range("L124").Value = 25000
$ws2.Range("N124").Value = -34820

# BSM row 105
$ws3.Range("H105").Value = 10198.857
$ws3.Range("I105").Value = 9378.4
$ws3.Range("K105").Value = 9378.4
$ws3.Range("M105").Value = -7631.4

# CRP row 31
$ws4.Range("H31").Value = 6252.674
$ws4.Range("I31").Value = 6047.925
$ws4.Range("J31").Value = 7617.6665
$ws4.Range("K31").Value = 6047.925
$ws4.Range("L31").Value = 7617.6665
$ws4.Range("M31").Value = -5752.925
$ws4.Range("N31").Value = -8207.666499999999

# CRP row 34
$ws4.Range("H34").Value = 6252.674
$ws4.Range("I34").Value = 6047.925
$ws4.Range("J34").Value = 7617.6665
$ws4.Range("K34").Value = 6047.925
$ws4.Range("L34").Value = 7617.6665
$ws4.Range("M34").Value = -5845.925
$ws4.Range("N34").Value = -8021.6665

# CRP row 58
$ws4.Range("H58").Value = 3762.4375
$ws4.Range("I58").Value = 3829.5334
$ws4.Range("J58").Value = 2756
$ws4.Range("K58").Value = 3829.5334
$ws4.Range("L58").Value = 2756
$ws4.Range("M58").Value = -3626.5334
$ws4.Range("N58").Value = -3162

# CRP row 132
$ws4.Range("H132").Value = 1619.125
$ws4.Range("I132").Value = 1559.0869
$ws4.Range("J132").Value = 3000
$ws4.Range("K132").Value = 4677.2607
$ws4.Range("L132").Value = 9000
$ws4.Range("M132").Value = -2147.2607
$ws4.Range("N132").Value = -14060

# CRP row 136
$ws4.Range("H136").Value = 3762.4375
$ws4.Range("I136").Value = 3829.5334
$ws4.Range("J136").Value = 2756
$ws4.Range("K136").Value = 11488.6002
$ws4.Range("L136").Value = 8268
$ws4.Range("M136").Value = -8938.600199999999
$ws4.Range("N136").Value = -13368

# CUL row 101
$ws5.Range("H101").Value = 0
$ws5.Range("J101").Value = 0
$ws5.Range("N101").ClearContents()

# CUL row 105
$ws5.Range("H105").Value = 0
$ws5.Range("J105").Value = 0
$ws5.Range("N105").ClearContents()

# CUL row 131
$ws5.Range("H131").Value = 1663.36
$ws5.Range("I131").Value = 946
$ws5.Range("J131").Value = 2440.5
$ws5.Range("K131").Value = 2838
$ws5.Range("L131").Value = 7321.5
$ws5.Range("M131").Value = 2202
$ws5.Range("N131").Value = -17401.5

# GSM row 107
$ws6.Range("H107").Value = 734
$ws6.Range("I107").Value = 593.7646999999999
$ws6.Range("J107").Value = 1528.6666
$ws6.Range("K107").Value = 593.7646999999999
$ws6.Range("L107").Value = 1528.6666
$ws6.Range("M107").Value = 1326.2353
$ws6.Range("N107").Value = -5368.6666

# GSM row 113
$ws6.Range("H113").Value = 614
$ws6.Range("I113").Value = 614
$ws6.Range("K113").Value = 614
$ws6.Range("M113").Value = 1556

# GSM row 123
$ws6.Range("H123").Value = 74997.5
$ws6.Range("J123").Value = 74997.5
$ws6.Range("L123").Value = 74997.5
$ws6.Range("N123").Value = -79897.5

# GSM row 126
$ws6.Range("H126").Value = 8819.1
$ws6.Range("I126").Value = 8399.5
$ws6.Range("K126").Value = 25198.5
$ws6.Range("M126").Value = -22728.5

# GSM row 132
$ws6.Range("H132").Value = 4886.533
$ws6.Range("I132").Value = 4967.577
$ws6.Range("J132").Value = 4359.75
$ws6.Range("K132").Value = 14902.731
$ws6.Range("L132").Value = 13079.25
$ws6.Range("M132").Value = -12372.731
$ws6.Range("N132").Value = -18139.25

# LTW row 132
$ws7.Range("H132").Value = 2106.8333
$ws7.Range("I132").Value = 2273.5
$ws7.Range("J132").Value = 1773.5
$ws7.Range("K132").Value = 6820.5
$ws7.Range("L132").Value = 5320.5
$ws7.Range("M132").Value = -4290.5
$ws7.Range("N132").Value = -10380.5

# WVR row 46
$ws8.Range("H46").Value = 66643.60000000001
$ws8.Range("J46").Value = 64304.5
$ws8.Range("L46").Value = 64304.5
$ws8.Range("N46").Value = -64766.5

# WVR row 62
$ws8.Range("H62").Value = 4699.6665
$ws8.Range("I62").Value = 3499.5
$ws8.Range("J62").Value = 5299.75
$ws8.Range("K62").Value = 3499.5
$ws8.Range("L62").Value = 5299.75
$ws8.Range("M62").Value = -2875.5
$ws8.Range("N62").Value = -6547.75

# WVR row 65
$ws8.Range("H65").Value = 4699.6665
$ws8.Range("I65").Value = 3499.5
$ws8.Range("J65").Value = 5299.75
$ws8.Range("K65").Value = 17497.5
$ws8.Range("L65").Value = 26498.75
$ws8.Range("M65").Value = -14377.5
$ws8.Range("N65").Value = -32738.75

# WVR row 81
$ws8.Range("H81").Value = 3232.8333
$ws8.Range("I81").Value = 4932.6665
$ws8.Range("J81").Value = 1533
$ws8.Range("K81").Value = 9865.333000000001
$ws8.Range("L81").Value = 3066
$ws8.Range("M81").Value = -8804.333000000001
$ws8.Range("N81").Value = -5188

# WVR row 84
$ws8.Range("H84").Value = 3232.8333
$ws8.Range("I84").Value = 4932.6665
$ws8.Range("J84").Value = 1533
$ws8.Range("K84").Value = 49326.665
$ws8.Range("L84").Value = 15330
$ws8.Range("M84").Value = -44022.665
$ws8.Range("N84").Value = -25938

# WVR row 111
$ws8.Range("H111").Value = 11644
$ws8.Range("J111").Value = 11644
$ws8.Range("L111").Value = 11644
$ws8.Range("N111").Value = -19824

# WVR row 134
$ws8.Range("H134").Value = 66643.60000000001
$ws8.Range("J134").Value = 64304.5
$ws8.Range("L134").Value = 192913.5
$ws8.Range("N134").Value = -197983.5
